$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.16645348072052
$ws.Range("B1").Value = 2.43521785736084
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.367369651794434
$ws.Range("E1").Value = 1.234577894210815
